$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Cyril Abtan"
$ws.Range("B2").Value = "Cyril Abtan"

$ws.Range("A3").Value = "Partial Cyril Abtan"
$ws.Range("B3").Value = "Cyri "

$ws.Range("A4").Value = "Surname Cyril Abtan"
$ws.Range("B4").Value = "Abtan"

$ws.Range("A5").Value = "Firstname Cyricl Abtan"
$ws.Range("B5").Value = "Cyril"

$ws.Range("A2").Select()
